$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Near the end of the document: remove the duplicated bold
#    "Play Alien Antix Free Online Slot Review" paragraph, and replace
#    the text of the following italic paragraph (formerly the meta
#    description) with the new image-generation prompt.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t -like "*Play Alien Antix Free Online Slot Review*") {
        if ($i -gt 1) {
            $p.Range.Delete()
        }
    }
}

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t -like "*Experience unique gameplay with Alien Antix Slot*") {
        $r = $p.Range
        $textRange = $d.Range($r.Start, $r.End - 1)
        $textRange.Text = "Create an eye-catching feature image for the online slot game ""Alien Antix"". The image should be in a cartoon style and feature a happy Maya warrior wearing glasses. Make sure to incorporate elements of space and aliens in the image to match the theme of the game. The image should be vibrant and colorful, with the Maya warrior as the central focus, surrounded by aliens and other space objects. Make the image stand out to attract potential players to the game."
    }
}

# ------------------------------------------------------------------
# 2) Right after the H1 title, insert a new paragraph holding the
#    "Meta description" label (bold) followed by its text.
# ------------------------------------------------------------------
$title = $d.Paragraphs(1)
$title.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"
$metaRange = $metaPara.Range
$metaRange.InsertAfter("Meta description: Experience unique gameplay with Alien Antix Slot, featuring bonuses and excellent graphics. Play for free and win big!")

$labelStart = $metaPara.Range.Start
$labelEnd = $labelStart + 16
$labelRange = $d.Range($labelStart, $labelEnd)
$labelRange.Bold = 1

Write-Output "done"
